$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add two new trailing columns: "date dep" (N) and "date arr" (O)
$ws.Range("N1").Value = "date dep"
$ws.Range("O1").Value = "date arr"

# Fill in the "date arr" values for the three data rows, matching the
# existing date column's dd/mm/yy display format. Setting NumberFormat
# before Value avoids Excel auto-creating a throwaway m/d/yyyy style.
$ws.Range("O2").NumberFormat = "dd/mm/yy"
$ws.Range("O2").Value = "5/22/2022"

$ws.Range("O3").NumberFormat = "dd/mm/yy"
$ws.Range("O3").Value = "1/1/2014"

$ws.Range("O4").NumberFormat = "dd/mm/yy"
$ws.Range("O4").Value = "1/1/2015"

# Move the active selection to N5, mirroring the original D5 selection
# shifted to sit below the newly added columns.
$ws.Range("N5").Select()
